$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.284.24'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '1.907.25'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '0.724'
$ws.Range('E5').Value = '  +9.03%  '
$ws.Range('D6').Value = '255.84'
$ws.Range('E6').Value = '  +3.95%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '40.71'
$ws.Range('E8').Value = '  -1.43%  '
$ws.Range('D9').Value = '0.378'
$ws.Range('E9').Value = '  +8.39%  '
$ws.Range('D10').Value = '52.85'
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('E11').Value = '  +5.34%  '
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').Value = '2.186.03'
$ws.Range('E13').Value = '  +0.35%  '
$ws.Range('D14').Value = '12.91'
$ws.Range('E14').Value = '  +6.59%  '
$ws.Range('D15').Value = '0.728'
$ws.Range('E15').Value = '  +4.47%  '
$ws.Range('D16').Value = '4.97'
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('D17').Value = '1.885.22'
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('D18').Value = '35.275.40'
$ws.Range('E18').Value = '  -0.19%  '
$ws.Range('D19').Value = '74.82'
$ws.Range('E19').Value = '  +3.59%  '
$ws.Range('D20').Value = '0.0₃0850'
$ws.Range('E20').Value = '  +3.52%  '
$ws.Range('D21').Value = '243.78'
$ws.Range('E21').Value = '  +1.34%  '
$ws.Range('D22').Value = '13.05'
$ws.Range('E22').Value = '  +4.81%  '
$ws.Range('D23').Value = '5.12'
$ws.Range('E23').Value = '  +5.71%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').Value = '2.46'
$ws.Range('E25').Value = '  +7.25%  '
$ws.Range('D26').Value = '2.45'
$ws.Range('E26').Value = '  +4.87%  '
$ws.Range('D27').Value = '166.26'
$ws.Range('E27').Value = '  -2.14%  '
$ws.Range('D28').Value = '8.70'
$ws.Range('E28').Value = '  +3.38%  '
$ws.Range('D29').Value = '18.74'
$ws.Range('E29').Value = '  +2.17%  '
$ws.Range('E30').Value = '  +4.44%  '
$ws.Range('D31').Value = '4.129.26'
$ws.Range('E32').Value = '  +6.38%  '
$ws.Range('E33').Value = '  +14.74%  '
$ws.Range('B34').Value = 'TrustWalletToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D34').Value = '1.64'
$ws.Range('E34').Value = '  +22.26%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.0589'
$ws.Range('E35').Value = '  +4.27%  '
$ws.Range('D36').Value = '4.26'
$ws.Range('E36').Value = '  +4.09%  '
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('D38').Value = '0.912'
$ws.Range('E38').Value = '  -2.32%  '
$ws.Range('E39').Value = '  +0.48%  '
$ws.Range('E40').Value = '  +5.27%  '
$ws.Range('D41').Value = '17.16'
$ws.Range('E41').Value = '  +6.47%  '
$ws.Range('D42').Value = '96.52'
$ws.Range('E42').Value = '  +7.74%  '
$ws.Range('E43').Value = '  +2.08%  '
$ws.Range('D44').Value = '0.0651'
$ws.Range('E44').Value = '  +4.04%  '
$ws.Range('D45').Value = '1.338.32'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').Value = '2.45'
$ws.Range('E46').Value = '  +2.68%  '
$ws.Range('D47').Value = '2.43'
$ws.Range('E47').Value = '  +0.93%  '
$ws.Range('D48').Value = '6.73'
$ws.Range('E48').Value = '  +3.59%  '
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('D50').Value = '45.17'
$ws.Range('E50').Value = '  -6.25%  '
$ws.Range('D51').Value = '0.0754'
$ws.Range('E51').Value = '  +6.95%  '
